$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "K, Pa"
$ws.Range("A3").Value = "n"
